# Update the "score" column (E) from a 0-5 scale to a 0-100 scale
# (multiply each existing non-zero score by 20) for the rows that changed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(2,4,5,6,8,9,10,12,13,26,34,38,40,41,42,44,45,46,48,49,50,52,53,54,55,56,57,58,60,61,62,63,64,65,70,72,73,74,75,76,77,78,79,80,81,82,83,84,85,98,99,110,111,112,113,114,115,116,117,118,119,120,121,122,123,124,125,126,127,128,129,130,131,132,133,134,135,136,137,138,139,140,141,142,143,144,145,146,147,148,149,150,151,152,153,154,155,156,157,159,163,166,167,170,175,183,184,185,186,187,188,189,195,196,197,198,199,200,201,202,203,204,205,206,207,208,209,211,212,213,215,216,217)
$newValues = @(100,60,20,40,60,80,40,60,20,20,100,100,100,100,100,20,60,100,40,20,100,100,100,100,100,100,100,100,100,100,100,100,100,100,100,40,20,100,100,60,100,100,100,100,80,100,100,100,100,60,60,100,100,100,100,100,100,100,100,100,100,60,60,100,100,100,100,100,100,100,100,100,100,100,100,100,100,100,100,100,100,60,60,100,100,100,100,100,100,100,100,100,100,100,80,100,100,100,80,40,20,100,20,100,20,100,100,100,100,100,100,100,100,100,100,100,100,100,100,100,100,100,60,20,100,100,100,100,100,100,100,20,100)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $ws.Cells.Item($rows[$i], 5).Value = $newValues[$i]
}
